$wb = $excel.ActiveWorkbook

# --- "user" sheet: fix "Password" -> "password" column label ---
$wsUser = $wb.Worksheets.Item("user")
$wsUser.Range("B1").Value = "password"

# --- "chatHistory" sheet: insert a new "type" column (E) describing the
#     kind of the "info" column (now shifted to F), and document the new
#     message/file/user-add/user-delete history rows ---
$wsChatHistory = $wb.Worksheets.Item("chatHistory")
$wsChatHistory.Columns("E:E").Insert()

$wsChatHistory.Range("E1").Value = "type"
$wsChatHistory.Range("E2").Value = "varchar(50)"
$wsChatHistory.Range("E3").Value = "not null"
$wsChatHistory.Range("E4").Value = "메시지 종류"
$wsChatHistory.Range("E5").Value = "메시지"

$wsChatHistory.Range("E6").Value = "파일"
$wsChatHistory.Range("F6").Value = "파일ID"

$wsChatHistory.Range("E7").Value = "사용자추가"
$wsChatHistory.Range("F7").Value = "user name"

$wsChatHistory.Range("E8").Value = "사용자삭제"
$wsChatHistory.Range("F8").Value = "user name"

[void]$wsChatHistory.Range("D5").Select()

# --- tab/selection bookkeeping: "user" becomes the active sheet/tab
#     (previously "chat" was active) ---
$wsUser.Activate()
[void]$wsUser.Range("B1").Select()
